$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 65

# --- Column A (Indice) : bold + bordered style like the other index cells ---
$ws.Range("A64").Copy()
$ws.Range("A65").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item($row, 1).Value = 64

# --- Column E (data_partida) : date/time number-format style like other rows ---
$ws.Range("E64").Copy()
$ws.Range("E65").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item($row, 5).Value = 45242.66666666666

# --- Remaining plain columns ---
$ws.Cells.Item($row, 2).Value = "azerbaijan"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 6).Value = "Qarabag"
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Value = "Sabail"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1.19
$ws.Cells.Item($row, 11).Value = "11/11/2023 04:13"
$ws.Cells.Item($row, 12).Value = 1.14
$ws.Cells.Item($row, 13).Value = "12/11/2023 15:47"
$ws.Cells.Item($row, 14).Value = 5.95
$ws.Cells.Item($row, 15).Value = "11/11/2023 04:13"
$ws.Cells.Item($row, 16).Value = 7.82
$ws.Cells.Item($row, 17).Value = "12/11/2023 15:47"
$ws.Cells.Item($row, 18).Value = 9.220000000000001
$ws.Cells.Item($row, 19).Value = "11/11/2023 04:13"
$ws.Cells.Item($row, 20).Value = 16.52
$ws.Cells.Item($row, 21).Value = "12/11/2023 15:47"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/qarabag-agdam-sabail/fikIS75A/"

Write-Output "row 65 written"
